# Insert one new weekly price record for "Zapallo italiano" (Macroferia
# Regional de Talca) directly above the existing row 536. Excel shifts the
# remaining historical rows (old 536-580) down by one, so the sheet grows
# from A1:R580 to A1:R581 - matching the dimension change in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 536..580 down to 537..581, leaving a blank row 536 to fill in.
$ws.Rows(536).Insert()

$ws.Range("A536").Value2 = 5
$ws.Range("B536").Value2 = "Macroferia Regional de Talca"
$ws.Range("C536").Value2 = "Maule"
$ws.Range("D536").Value2 = 45106
$ws.Range("E536").Value2 = 7
$ws.Range("F536").Value2 = 100112032
$ws.Range("G536").Value2 = "Zapallo italiano"
$ws.Range("H536").Value2 = "Sin especificar"
$ws.Range("I536").Value2 = "Primera"
$ws.Range("J536").Value2 = 300
$ws.Range("K536").Value2 = 15000
$ws.Range("L536").Value2 = 15000
$ws.Range("M536").Value2 = 15000
$ws.Range("N536").Value2 = "`$/caja 50 unidades"
$ws.Range("O536").Value2 = "Región de Arica y Parinacota"
$ws.Range("P536").Value2 = 300
$ws.Range("Q536").Value2 = 50
$ws.Range("R536").Value2 = "Hortaliza"
